# cambodia-file.xlsx edit: add a "Type" column and format the phone-number
# column as Text so a leading "+" can be entered without Excel mangling it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet only ever has a header row; rows 2-10 are empty placeholder rows
# left over from the template. Drop them so the used range collapses back
# down to just the header row (dimension goes from A1:O10 to A1:P1).
$ws.Rows("2:10").Delete()

# New header in column P: "Type"
$ws.Cells.Item(1, 16).Value = "Type"

# Column I is "Beneficiary Phone Number" - format it as Text (the "@"
# number format) so phone numbers can be entered with a leading "+"
# without Excel coercing/clipping them as numeric values.
$ws.Cells.Item(1, 9).NumberFormat = "@"

# Leave the selection where it was when the file was last saved.
[void]$ws.Range("H8").Select()
